$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

# Row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# Row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 1505.6666
$ws.Range("I92").Value = 1099.7142
$ws.Range("K92").Value = 1099.7142
$ws.Range("M92").Value = 148.2858000000001

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1599.5
$ws.Range("I98").Value = 1199
$ws.Range("K98").Value = 1199
$ws.Range("M98").Value = 299

# Row 104 (Leve Item ID 24263)
$ws.Range("H104").Value = 8300
$ws.Range("I104").Value = 8300
$ws.Range("K104").Value = 24900
$ws.Range("M104").Value = -23153

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1599.5
$ws.Range("I122").Value = 1199
$ws.Range("K122").Value = 3597
$ws.Range("M122").Value = -1147

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2666.1667
$ws.Range("I138").Value = 997
$ws.Range("K138").Value = 2991
$ws.Range("M138").Value = 2149

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 6383.476
$ws.Range("I32").Value = 3558.5
$ws.Range("K32").Value = 3558.5
$ws.Range("M32").Value = -3271.5

# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 1496.3334
$ws.Range("I45").Value = 1496.3334
$ws.Range("K45").Value = 1496.3334
$ws.Range("M45").Value = -1119.3334

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9788

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 4600
$ws.Range("I122").Value = 3966.6667
$ws.Range("K122").Value = 11900.0001
$ws.Range("M122").Value = -9450.000100000001

# Row 131 (Leve Item ID 34706)
$ws.Range("H131").Value = 75000
$ws.Range("J131").Value = 75000
$ws.Range("L131").Value = 75000
$ws.Range("N131").Value = -85080

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 10000
$ws.Range("I136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("M136").Value = -27450

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9753

$ws = $wb.Worksheets.Item("CRP")
# Row 53 (Leve Item ID 25632)
$ws.Range("H53").Value = 43142
$ws.Range("J53").Value = 43142
$ws.Range("L53").Value = 43142
$ws.Range("N53").Value = -44356

# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 5049.25
$ws.Range("I62").Value = 4998.5
$ws.Range("K62").Value = 4998.5
$ws.Range("M62").Value = -4374.5

# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 5049.25
$ws.Range("I65").Value = 4998.5
$ws.Range("K65").Value = 24992.5
$ws.Range("M65").Value = -21872.5

# Row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 94 (Leve Item ID 32934)
$ws.Range("H94").Value = 4070
$ws.Range("I94").Value = 5783.3335
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 5783.3335
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = -5332.3335
$ws.Range("N94").Value = -2402

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1275
$ws.Range("I122").Value = 1362.5
$ws.Range("K122").Value = 4087.5
$ws.Range("M122").Value = -1637.5

# Row 141 (Leve Item ID 43345)
$ws.Range("H141").Value = 795776.7
$ws.Range("J141").Value = 795776.7
$ws.Range("L141").Value = 795776.7
$ws.Range("N141").Value = -806136.7

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 60.2
$ws.Range("I2").Value = 53
$ws.Range("J2").Value = 65
$ws.Range("K2").Value = 318
$ws.Range("L2").Value = 390
$ws.Range("M2").Value = -205
$ws.Range("N2").Value = -616

# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1060
$ws.Range("I68").Value = 800
$ws.Range("J68").Value = 1125
$ws.Range("K68").Value = 2400
$ws.Range("L68").Value = 3375
$ws.Range("M68").Value = -1589
$ws.Range("N68").Value = -4997

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1060
$ws.Range("I71").Value = 800
$ws.Range("J71").Value = 1125
$ws.Range("K71").Value = 7200
$ws.Range("L71").Value = 10125
$ws.Range("M71").Value = -3144
$ws.Range("N71").Value = -18237

# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 6957.143
$ws.Range("I70").Value = 6926.25
$ws.Range("K70").Value = 6926.25
$ws.Range("M70").Value = -6656.25

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 6957.143
$ws.Range("I73").Value = 6926.25
$ws.Range("K73").Value = 6926.25
$ws.Range("M73").Value = -5990.25

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 933.3333
$ws.Range("I16").Value = 900.5
$ws.Range("J16").Value = 999
$ws.Range("K16").Value = 900.5
$ws.Range("L16").Value = 999
$ws.Range("M16").Value = -730.5
$ws.Range("N16").Value = -1339

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 865.38464
$ws.Range("J46").Value = 937.5
$ws.Range("L46").Value = 937.5
$ws.Range("N46").Value = -1313.5

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 599.5
$ws.Range("I55").Value = 600
$ws.Range("K55").Value = 600
$ws.Range("M55").Value = -427

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 7748.75
$ws.Range("I122").Value = 4331.6665
$ws.Range("K122").Value = 12994.9995
$ws.Range("M122").Value = -10544.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 14 (Leve Item ID 2658)
$ws.Range("H14").Value = 1442.0834
$ws.Range("I14").Value = 1442.0834
$ws.Range("K14").Value = 1442.0834
$ws.Range("M14").Value = -1274.0834

# Row 18 (Leve Item ID 3543)
$ws.Range("H18").Value = 99998.5
$ws.Range("J18").Value = 99998.5
$ws.Range("L18").Value = 99998.5
$ws.Range("N18").Value = -100344.5

# Row 30 (Leve Item ID 2700)
$ws.Range("H30").Value = 5000
$ws.Range("I30").Value = 5000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 5000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -4893
$ws.Range("N30").ClearContents()

# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 451
$ws.Range("I113").Value = 451
$ws.Range("K113").Value = 1353
$ws.Range("M113").Value = 817

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 15521.286
$ws.Range("I122").Value = 1383.3334
$ws.Range("K122").Value = 4150.0002
$ws.Range("M122").Value = -1700.0002

# Row 123 (Leve Item ID 34127)
$ws.Range("H123").Value = 158000
$ws.Range("J123").Value = 158000
$ws.Range("L123").Value = 158000
$ws.Range("N123").Value = -167800

# Row 130 (Leve Item ID 34705)
$ws.Range("H130").Value = 30285.666
$ws.Range("J130").Value = 30285.666
$ws.Range("L130").Value = 30285.666
$ws.Range("N130").Value = -40325.666

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 3993.75
$ws.Range("I136").Value = 3993.75
$ws.Range("K136").Value = 11981.25
$ws.Range("M136").Value = -9431.25
